$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2305
$ws.Range("I70").Value = 1151.25
$ws.Range("J70").Value = 2881.875
$ws.Range("K70").Value = 3453.75
$ws.Range("L70").Value = 8645.625
$ws.Range("M70").Value = -3183.75
$ws.Range("N70").Value = -9185.625
$ws.Range("H73").Value = 2305
$ws.Range("I73").Value = 1151.25
$ws.Range("J73").Value = 2881.875
$ws.Range("K73").Value = 3453.75
$ws.Range("L73").Value = 8645.625
$ws.Range("M73").Value = -2517.75
$ws.Range("N73").Value = -10517.625
$ws.Range("H86").Value = 1354.3077
$ws.Range("I86").Value = 1289.5
$ws.Range("J86").Value = 1570.3334
$ws.Range("K86").Value = 1289.5
$ws.Range("L86").Value = 1570.3334
$ws.Range("M86").Value = -166.5
$ws.Range("N86").Value = -3816.3334
$ws.Range("H89").Value = 1354.3077
$ws.Range("I89").Value = 1289.5
$ws.Range("J89").Value = 1570.3334
$ws.Range("K89").Value = 6447.5
$ws.Range("L89").Value = 7851.666999999999
$ws.Range("M89").Value = -831.5
$ws.Range("N89").Value = -19083.667
$ws.Range("H106").Value = 7008544.5
$ws.Range("I106").Value = 8625371
$ws.Range("J106").Value = 2298.6667
$ws.Range("K106").Value = 8625371
$ws.Range("L106").Value = 2298.6667
$ws.Range("M106").Value = -8624740
$ws.Range("N106").Value = -3560.6667

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20530.29
$ws.Range("I32").Value = 1955.8334
$ws.Range("J32").Value = 66420.12
$ws.Range("K32").Value = 1955.8334
$ws.Range("L32").Value = 66420.12
$ws.Range("M32").Value = -1668.8334
$ws.Range("N32").Value = -66994.12
$ws.Range("H88").Value = 6401.5
$ws.Range("I88").Value = 1603
$ws.Range("J88").Value = 11200
$ws.Range("K88").Value = 1603
$ws.Range("L88").Value = 11200
$ws.Range("M88").Value = -1197
$ws.Range("N88").Value = -12012
$ws.Range("H91").Value = 6401.5
$ws.Range("I91").Value = 1603
$ws.Range("J91").Value = 11200
$ws.Range("K91").Value = 1603
$ws.Range("L91").Value = 11200
$ws.Range("M91").Value = -199
$ws.Range("N91").Value = -14008

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5692.9165
$ws.Range("I86").Value = 1489.3889
$ws.Range("J86").Value = 18303.5
$ws.Range("K86").Value = 1489.3889
$ws.Range("L86").Value = 18303.5
$ws.Range("M86").Value = -366.3888999999999
$ws.Range("N86").Value = -20549.5
$ws.Range("H89").Value = 5692.9165
$ws.Range("I89").Value = 1489.3889
$ws.Range("J89").Value = 18303.5
$ws.Range("K89").Value = 7446.9445
$ws.Range("L89").Value = 91517.5
$ws.Range("M89").Value = -1830.9445
$ws.Range("N89").Value = -102749.5
$ws.Range("H94").Value = 718.3333
$ws.Range("I94").Value = 702
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 702
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -251
$ws.Range("N94").Value = -1702
$ws.Range("H107").Value = 293.80768
$ws.Range("I107").Value = 169.06667
$ws.Range("J107").Value = 463.9091
$ws.Range("K107").Value = 169.06667
$ws.Range("L107").Value = 463.9091
$ws.Range("M107").Value = 1750.93333
$ws.Range("N107").Value = -4303.9091

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 24700
$ws.Range("I62").Value = 31428.572
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 31428.572
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -30804.572
$ws.Range("N62").Value = -10248
$ws.Range("H65").Value = 24700
$ws.Range("I65").Value = 31428.572
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 157142.86
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -154022.86
$ws.Range("N65").Value = -51240

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 15631.5
$ws.Range("I68").Value = 40267.332
$ws.Range("J68").Value = 850
$ws.Range("K68").Value = 120801.996
$ws.Range("L68").Value = 2550
$ws.Range("M68").Value = -119990.996
$ws.Range("N68").Value = -4172
$ws.Range("H71").Value = 15631.5
$ws.Range("I71").Value = 40267.332
$ws.Range("J71").Value = 850
$ws.Range("K71").Value = 362405.988
$ws.Range("L71").Value = 7650
$ws.Range("M71").Value = -358349.988
$ws.Range("N71").Value = -15762
$ws.Range("H87").Value = 12249.75
$ws.Range("I87").Value = 6599.6
$ws.Range("J87").Value = 21666.666
$ws.Range("K87").Value = 19798.8
$ws.Range("L87").Value = 64999.99800000001
$ws.Range("M87").Value = -18550.8
$ws.Range("N87").Value = -67495.99800000001
$ws.Range("H90").Value = 12249.75
$ws.Range("I90").Value = 6599.6
$ws.Range("J90").Value = 21666.666
$ws.Range("K90").Value = 59396.4
$ws.Range("L90").Value = 194999.994
$ws.Range("M90").Value = -53156.4
$ws.Range("N90").Value = -207479.994
$ws.Range("H109").Value = 1744.7778
$ws.Range("I109").Value = 175.75
$ws.Range("J109").Value = 3000
$ws.Range("K109").Value = 527.25
$ws.Range("L109").Value = 9000
$ws.Range("M109").Value = 512.75
$ws.Range("N109").Value = -11080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11196.667
$ws.Range("I80").Value = 12736
$ws.Range("K80").Value = 12736
$ws.Range("M80").Value = -11738
$ws.Range("H83").Value = 11196.667
$ws.Range("I83").Value = 12736
$ws.Range("K83").Value = 63680
$ws.Range("M83").Value = -58688
$ws.Range("H93").Value = 29164
$ws.Range("J93").Value = 29164
$ws.Range("L93").Value = 29164
$ws.Range("N93").Value = -32908
$ws.Range("H107").Value = 604.0526
$ws.Range("I107").Value = 979.6
$ws.Range("J107").Value = 186.77777
$ws.Range("K107").Value = 979.6
$ws.Range("L107").Value = 186.77777
$ws.Range("M107").Value = 940.4
$ws.Range("N107").Value = -4026.77777
$ws.Range("H113").Value = 1986.1428
$ws.Range("I113").Value = 1958
$ws.Range("J113").Value = 2014.2858
$ws.Range("K113").Value = 1958
$ws.Range("L113").Value = 2014.2858
$ws.Range("M113").Value = 212
$ws.Range("N113").Value = -6354.2858

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 368.5
$ws.Range("I46").Value = 422.8
$ws.Range("J46").Value = 314.2
$ws.Range("K46").Value = 422.8
$ws.Range("L46").Value = 314.2
$ws.Range("M46").Value = -234.8
$ws.Range("N46").Value = -690.2
$ws.Range("H93").Value = 533.5517
$ws.Range("I93").Value = 500.16666
$ws.Range("J93").Value = 693.8
$ws.Range("K93").Value = 500.16666
$ws.Range("L93").Value = 693.8
$ws.Range("M93").Value = 747.83334
$ws.Range("N93").Value = -3189.8
$ws.Range("H122").Value = 3053.1924
$ws.Range("I122").Value = 1923.8182
$ws.Range("K122").Value = 5771.4546
$ws.Range("M122").Value = -3321.4546
